$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 02:16"

# Swap row 47 (Lugo) and row 48 (Almeria) entirely, which also
# reflects the shared-string reorder (Almeria now precedes Lugo)
# and the swapped "Casos activos" (column C) values.
$a47 = $ws.Range("A47").Value2
$b47 = $ws.Range("B47").Value2
$c47 = $ws.Range("C47").Value2
$d47 = $ws.Range("D47").Value2
$e47 = $ws.Range("E47").Value2

$a48 = $ws.Range("A48").Value2
$b48 = $ws.Range("B48").Value2
$c48 = $ws.Range("C48").Value2
$d48 = $ws.Range("D48").Value2
$e48 = $ws.Range("E48").Value2

$ws.Range("A47").Value = $a48
$ws.Range("B47").Value = $b48
$ws.Range("C47").Value = $c48
$ws.Range("D47").Value = $d48
$ws.Range("E47").Value = $e48

$ws.Range("A48").Value = $a47
$ws.Range("B48").Value = $b47
$ws.Range("C48").Value = $c47
$ws.Range("D48").Value = $d47
$ws.Range("E48").Value = $e47
